$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "a=1,x=3"
$ws.Range("C2").Value = "1,2,3"
$ws.Range("I4").Value = "1,2,3,4"

$ws.Range("H3").Select()
